# Correção nos dados e início da análise PNAD 2009
#
# The "seguranca" worksheet had a data-entry offset bug: two rows that only
# carried a row label ("situação do domicílio" and "grandes regiões e
# unidades da federação") had no numeric data of their own, while every
# data row below them actually held the values that belonged to the row
# above it. Removing those two "empty" label rows shifts all the real
# figures up into their correctly-labelled rows (brasil / urbana / rural /
# norte / rondônia / ... / mato grosso / goiás / distrito federal) and
# naturally drops the two now-superfluous rows at the bottom of the table
# (the sheet shrinks from A1:I40 to A1:I38).
#
# Also, the column-B sub-header in row 2 ("unnamed: 1_level_1") is
# relabelled to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled sub-header in row 2.
$ws.Range("B2").Value = "total"

# Remove the row that only contains the "situação do domicílio" label
# (row 5). This shifts every row below it up by one, so the values that
# used to sit in row 6 now correctly line up with the "urbana" label, etc.
$ws.Range("A5").EntireRow.Delete()

# After the previous deletion, the row that only contains the "grandes
# regiões e unidades da federação" label has moved up to row 7. Remove it
# too, shifting everything below it up by one more row.
$ws.Range("A7").EntireRow.Delete()
